# Add three new trailing columns (D, E, F) with header labels to Sheet1,
# mirroring the upload that appended ORG_CLUB_IDENOLD / ORG_CLUB_IDENNEW /
# ORG_CLUB_STATUS to the ORG_CLUB export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ORG_CLUB_IDENOLD"
$ws.Range("E1").Value = "ORG_CLUB_IDENNEW"
$ws.Range("F1").Value = "ORG_CLUB_STATUS"

# Leave the selection where the author's save left it.
$ws.Range("F8").Select()
